# Fill in the second weekly-plan table (rows 12-16 under the
# "日期：2018.10.08 第六周周一" header) with each member's task and status,
# mirroring the layout of the first table (rows 3-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 - 王伟锋
$ws.Range("A12").Value = "王伟锋"
$ws.Range("B12").Value = "完成系统管理用例图"
$ws.Range("C12").Value = 1
$ws.Range("C12").NumberFormat = "0%"

# Row 13 - 陈升云
$ws.Range("A13").Value = "陈升云"
$ws.Range("B13").Value = "完善用户用例图"
$ws.Range("C13").Value = 1
$ws.Range("C13").NumberFormat = "0%"

# Row 14 - 林玮成
$ws.Range("A14").Value = "林玮成"
$ws.Range("B14").Value = "对用例图进行用例简述"
$ws.Range("C14").Value = "进行中"

# Row 15 - 吴帅辰
$ws.Range("A15").Value = "吴帅辰"
$ws.Range("B15").Value = "对用例图进行用例简述"
$ws.Range("C15").Value = "进行中"

# Row 16 - 李海洋
$ws.Range("A16").Value = "李海洋"
$ws.Range("B16").Value = "对用例图进行用例简述"
$ws.Range("C16").Value = "进行中"

# Match the author's final selection (cell C14 was last touched / selected).
$ws.Range("C14").Select() | Out-Null
